$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# Helper: find the 1-based index of the first paragraph whose trimmed text
# equals / starts-with a given string.
# ---------------------------------------------------------------------------
function Find-ParaIndexExact($doc, $text) {
    for ($i = 1; $i -le $doc.Paragraphs.Count; $i++) {
        $t = $doc.Paragraphs.Item($i).Range.Text
        $t = $t.TrimEnd([char]13, [char]7)
        if ($t -eq $text) {
            return $i
        }
    }
    return -1
}

function Find-ParaIndexStartsWith($doc, $text) {
    for ($i = 1; $i -le $doc.Paragraphs.Count; $i++) {
        $t = $doc.Paragraphs.Item($i).Range.Text
        if ($t.StartsWith($text)) {
            return $i
        }
    }
    return -1
}

# ---------------------------------------------------------------------------
# 1) Move the "_GoBack" bookmark away from its old spot (the blank paragraph
#    right before the first "Jeffrey" Code-Review-Summary heading).
# ---------------------------------------------------------------------------
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks.Item("_GoBack").Delete()
}

# ---------------------------------------------------------------------------
# 2) Six blank paragraph marks get inserted above the "Zifan" heading in the
#    Code Review Strategy section; the sz/szCs run-properties that used to
#    sit on that heading's paragraph mark move down onto the final (7th)
#    paragraph, which also now carries the relocated "_GoBack" bookmark.
# ---------------------------------------------------------------------------
$zifanIdx = Find-ParaIndexExact $d "Zifan"
$r = $d.Paragraphs.Item($zifanIdx).Range
$frag = @'
<?xml version="1.0" standalone="yes"?>
<?mso-application progid="Word.Document"?>
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">
<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">
<pkg:xmlData>
<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
<w:body>
<w:p>
  <w:pPr>
    <w:pStyle w:val="Heading2"/>
    <w:rPr>
      <w:rFonts w:ascii="Arial" w:eastAsia="Times New Roman" w:hAnsi="Arial" w:cs="Arial"/>
      <w:lang w:eastAsia="en-CA"/>
    </w:rPr>
  </w:pPr>
</w:p>
<w:p>
  <w:pPr>
    <w:pStyle w:val="Heading2"/>
    <w:rPr>
      <w:rFonts w:ascii="Arial" w:eastAsia="Times New Roman" w:hAnsi="Arial" w:cs="Arial"/>
      <w:lang w:eastAsia="en-CA"/>
    </w:rPr>
  </w:pPr>
</w:p>
<w:p>
  <w:pPr>
    <w:pStyle w:val="Heading2"/>
    <w:rPr>
      <w:rFonts w:ascii="Arial" w:eastAsia="Times New Roman" w:hAnsi="Arial" w:cs="Arial"/>
      <w:lang w:eastAsia="en-CA"/>
    </w:rPr>
  </w:pPr>
</w:p>
<w:p>
  <w:pPr>
    <w:pStyle w:val="Heading2"/>
    <w:rPr>
      <w:rFonts w:ascii="Arial" w:eastAsia="Times New Roman" w:hAnsi="Arial" w:cs="Arial"/>
      <w:lang w:eastAsia="en-CA"/>
    </w:rPr>
  </w:pPr>
</w:p>
<w:p/>
<w:p>
  <w:pPr>
    <w:pStyle w:val="Heading2"/>
    <w:rPr>
      <w:rFonts w:ascii="Arial" w:eastAsia="Times New Roman" w:hAnsi="Arial" w:cs="Arial"/>
      <w:lang w:eastAsia="en-CA"/>
    </w:rPr>
  </w:pPr>
</w:p>
<w:p>
  <w:pPr>
    <w:pStyle w:val="Heading2"/>
    <w:rPr>
      <w:rFonts w:ascii="Arial" w:eastAsia="Times New Roman" w:hAnsi="Arial" w:cs="Arial"/>
      <w:sz w:val="24"/>
      <w:szCs w:val="24"/>
      <w:lang w:eastAsia="en-CA"/>
    </w:rPr>
  </w:pPr>
  <w:bookmarkStart w:id="0" w:name="_GoBack"/>
  <w:bookmarkEnd w:id="0"/>
  <w:r>
    <w:rPr>
      <w:rFonts w:ascii="Arial" w:eastAsia="Times New Roman" w:hAnsi="Arial" w:cs="Arial"/>
      <w:lang w:eastAsia="en-CA"/>
    </w:rPr>
    <w:t>Zifan</w:t>
  </w:r>
</w:p>
</w:body>
</w:document>
</pkg:xmlData>
</pkg:part>
</pkg:package>
'@
$r.InsertXML($frag)

# ---------------------------------------------------------------------------
# 3) Merge the two runs that were split by a lastRenderedPageBreak in the
#    "DatabasePresetQuery.java effectively ... overall this part ..."
#    paragraph into a single run (drop the page-break artifact).
# ---------------------------------------------------------------------------
$idx = Find-ParaIndexStartsWith $d "The code is well formatted, meeting indentation"
$r = $d.Paragraphs.Item($idx).Range
$frag = @'
<?xml version="1.0" standalone="yes"?>
<?mso-application progid="Word.Document"?>
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">
<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">
<pkg:xmlData>
<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
<w:body>
<w:p>
  <w:pPr>
    <w:spacing w:after="0" w:line="240" w:lineRule="auto"/>
    <w:rPr>
      <w:rFonts w:ascii="-webkit-standard" w:hAnsi="-webkit-standard" w:cs="Times New Roman"/>
      <w:color w:val="000000"/>
      <w:lang w:val="en-US"/>
    </w:rPr>
  </w:pPr>
  <w:r>
    <w:rPr>
      <w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/>
      <w:color w:val="000000"/>
      <w:lang w:val="en-US"/>
    </w:rPr>
    <w:t xml:space="preserve">The code is well formatted, meeting indentation and whitespace standards and free from parse errors. The code did what is supposed to do. It includes helper functions to help the implementations in </w:t>
  </w:r>
  <w:r>
    <w:rPr>
      <w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/>
      <w:i/>
      <w:iCs/>
      <w:color w:val="000000"/>
      <w:lang w:val="en-US"/>
    </w:rPr>
    <w:t xml:space="preserve">DatabasePresetQuery.java </w:t>
  </w:r>
  <w:r>
    <w:rPr>
      <w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/>
      <w:color w:val="000000"/>
      <w:lang w:val="en-US"/>
    </w:rPr>
    <w:t xml:space="preserve">effectively. These helper functions are necessarily included in this class. It is nicely commented and can easily understand. I think overall this part of code is good because they are correct and effective solutions for the project requirements at hand. </w:t>
  </w:r>
</w:p>
</w:body>
</w:document>
</pkg:xmlData>
</pkg:part>
</pkg:package>
'@
$r.InsertXML($frag)

# ---------------------------------------------------------------------------
# 4) Add a lastRenderedPageBreak before the run of the first "Kelvin"
#    Code-Review-Summary heading.
# ---------------------------------------------------------------------------
$kelvinIdx = Find-ParaIndexExact $d "Kelvin"
$r = $d.Paragraphs.Item($kelvinIdx).Range
$frag = @'
<?xml version="1.0" standalone="yes"?>
<?mso-application progid="Word.Document"?>
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">
<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">
<pkg:xmlData>
<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
<w:body>
<w:p>
  <w:pPr>
    <w:pStyle w:val="Heading2"/>
    <w:rPr>
      <w:rFonts w:ascii="Arial" w:eastAsia="Times New Roman" w:hAnsi="Arial" w:cs="Arial"/>
      <w:sz w:val="24"/>
      <w:szCs w:val="24"/>
      <w:lang w:eastAsia="en-CA"/>
    </w:rPr>
  </w:pPr>
  <w:r>
    <w:rPr>
      <w:rFonts w:ascii="Arial" w:eastAsia="Times New Roman" w:hAnsi="Arial" w:cs="Arial"/>
      <w:lang w:eastAsia="en-CA"/>
    </w:rPr>
    <w:lastRenderedPageBreak/>
    <w:t>Kelvin</w:t>
  </w:r>
</w:p>
</w:body>
</w:document>
</pkg:xmlData>
</pkg:part>
</pkg:package>
'@
$r.InsertXML($frag)

# ---------------------------------------------------------------------------
# 5) Merge the two runs that were split by a lastRenderedPageBreak in the
#    "... functions are small (less than 10 lines) ..." paragraph into a
#    single run (drop the page-break artifact).
# ---------------------------------------------------------------------------
$idx = Find-ParaIndexStartsWith $d "The code is easily understandable, performs as expected and commented well."
$r = $d.Paragraphs.Item($idx).Range
$frag = @'
<?xml version="1.0" standalone="yes"?>
<?mso-application progid="Word.Document"?>
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">
<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">
<pkg:xmlData>
<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
<w:body>
<w:p>
  <w:pPr>
    <w:spacing w:after="0" w:line="240" w:lineRule="auto"/>
    <w:rPr>
      <w:rFonts w:ascii="Arial" w:eastAsia="Times New Roman" w:hAnsi="Arial" w:cs="Arial"/>
      <w:color w:val="000000"/>
      <w:lang w:val="en-US" w:eastAsia="en-CA"/>
    </w:rPr>
  </w:pPr>
  <w:r>
    <w:rPr>
      <w:rFonts w:ascii="Arial" w:eastAsia="Times New Roman" w:hAnsi="Arial" w:cs="Arial"/>
      <w:color w:val="000000"/>
      <w:lang w:val="en-US" w:eastAsia="en-CA"/>
    </w:rPr>
    <w:t>The code is easily understandable, performs as expected and commented well. The names of the functions explain what they do. Exceptions are handled properly by a try and catch statement. Docstrings explain the function of the code and functions are small (less than 10 lines). Good use of mock objects for testing. No dead code other than one unused import statement in TestService.java.</w:t>
  </w:r>
  <w:r>
    <w:rPr>
      <w:rFonts w:ascii="Arial" w:eastAsia="Times New Roman" w:hAnsi="Arial" w:cs="Arial"/>
      <w:i/>
      <w:iCs/>
      <w:color w:val="000000"/>
      <w:lang w:val="en-US" w:eastAsia="en-CA"/>
    </w:rPr>
    <w:t xml:space="preserve"> </w:t>
  </w:r>
  <w:r>
    <w:rPr>
      <w:rFonts w:ascii="Arial" w:eastAsia="Times New Roman" w:hAnsi="Arial" w:cs="Arial"/>
      <w:color w:val="000000"/>
      <w:lang w:val="en-US" w:eastAsia="en-CA"/>
    </w:rPr>
    <w:t>No copy and pasted code, and spacing is consistent throughout. Follows SOLID principles where a class only performs a single function. Performs as required according to our original design.</w:t>
  </w:r>
</w:p>
</w:body>
</w:document>
</pkg:xmlData>
</pkg:part>
</pkg:package>
'@
$r.InsertXML($frag)

Write-Output "done"
